# Update column C (Fitness) values on the active sheet.
# Rows 2-252 correspond to Generation 0-250 (column B), and the Fitness
# value (column C) is updated in five contiguous blocks as follows:
#   rows 2-19   (Generation 0-17)   : 7569 -> 7345
#   rows 20-24  (Generation 18-22)  : 7569 -> 7343
#   rows 25-110 (Generation 23-108) : 7569 -> 7310
#   rows 111-154(Generation 109-152): 7569 -> 7295
#   rows 155-252(Generation 153-250): 7569 -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$segments = @(
    @{ Start = 2;   End = 19;  Value = 7345 },
    @{ Start = 20;  End = 24;  Value = 7343 },
    @{ Start = 25;  End = 110; Value = 7310 },
    @{ Start = 111; End = 154; Value = 7295 },
    @{ Start = 155; End = 252; Value = 7293 }
)

foreach ($segment in $segments) {
    $rangeAddress = "C$($segment.Start):C$($segment.End)"
    $ws.Range($rangeAddress).Value = $segment.Value
}
